# Update cryptocurrency price/volume data (and swap two rows whose
# ranking order changed) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.393.49'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.821.92'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.17'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5237'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.24%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3851'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08031'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.77%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.116'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.90'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.395'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.91'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.420'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('D16').Value = '1.822.26'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '94.45'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001103'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.66'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('D23').Value = '28.450.49'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.39'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.03%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.242'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '159.31'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.87'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.54%  '
$ws.Range('D28').Value = '2.028.25'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.418'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.49'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1105'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +2.03%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.081'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.679'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('E34').Value = '  +0.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07362'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.93%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '12.32'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +9.92%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2199'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02342'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.149'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.738'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.26%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6329'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.183'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.380'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.50'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6139'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +4.17%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.785'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '127.33'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.23%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.983'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.207'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06898'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '73.78'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.75%  '
